$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row before row 10 ("Id" 111, Tutorial01/Next), shifting it and
# every row below down by one. This makes room for a new "Add Pawn" example
# row (Id 108) right after Id 107.
$ws.Rows.Item(10).Insert()

# Fill the new row 10 with the new "Add Pawn" example option data.
# (Inheriting the format/style of the surrounding data rows, as Excel does
# automatically on row insert.)
$ws.Range("A10").Value = 108
$ws.Range("B10").Value = "Option8"
$ws.Range("C10").Value = "加一个角色"
$ws.Range("D10").Value = "[[109:1003:26]]"

# Leave the selection on the newly-added row, matching the author's final
# cursor position.
$ws.Range("D10").Select()
